# feat: add tileset slicer and tileset-backed neon map
#
# Adds 5 new columns (J:K:L:M:N) to the "__data" sheet describing a
# tileset-backed slicer for the map tiles (tileSheet/tileWidth/tileHeight/
# tileRow/tileCol), mirroring the existing A:I schema (row 4 = type row,
# row 5 = field-name row, rows 6-9 = data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

# Row 4 - column type markers
$ws.Range("J4").Value = "'string"
$ws.Range("K4").Value = "'uint"
$ws.Range("L4").Value = "'uint"
$ws.Range("M4").Value = "'uint"
$ws.Range("N4").Value = "'uint"

# Row 5 - column names
$ws.Range("J5").Value = "'tileSheet"
$ws.Range("K5").Value = "'tileWidth"
$ws.Range("L5").Value = "'tileHeight"
$ws.Range("M5").Value = "'tileRow"
$ws.Range("N5").Value = "'tileCol"

# Row 6 - Ruined Arcade
$ws.Range("J6").Value = "'tileset_city"
$ws.Range("K6").Value = "'32"
$ws.Range("L6").Value = "'32"
$ws.Range("M6").Value = "'0"
$ws.Range("N6").Value = "'0"

# Row 7 - Subway Artery
$ws.Range("J7").Value = "'tileset_city"
$ws.Range("K7").Value = "'32"
$ws.Range("L7").Value = "'32"
$ws.Range("M7").Value = "'0"
$ws.Range("N7").Value = "'1"

# Row 8 - Harbor Breach
$ws.Range("J8").Value = "'tileset_city"
$ws.Range("K8").Value = "'32"
$ws.Range("L8").Value = "'32"
$ws.Range("M8").Value = "'1"
$ws.Range("N8").Value = "'0"

# Row 9 - Polar Lab Perimeter
$ws.Range("J9").Value = "'tileset_city"
$ws.Range("K9").Value = "'32"
$ws.Range("L9").Value = "'32"
$ws.Range("M9").Value = "'1"
$ws.Range("N9").Value = "'1"

# Keep the "numbers stored as text" error-check suppressed across the
# now-wider used range (A4:N9), matching the sheet's existing A4:I9 rule.
$fullRange = $ws.Range("A4:N9")
$fullRange.Errors.Item(9).Ignore = $true
